# Dotaznik workbook update: add "var" (variance) row and relabel rows with
# G-column captions (prum / var / vyb prum), matching the commit
# "Přidány úpravy podle Albertova".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new blank row 9; the old row 9 (selected-rows average) shifts to
# row 10, its shared formula splits into independent per-cell formulas.
$ws.Rows.Item(9).Insert()

# Row/column captions in column G. Order matters: new shared strings are
# appended in the order first written, and need to land as
# prum=72, vyb prum=73, var=74.
$ws.Range("G7").Value = "prum"
$ws.Range("G10").Value = "vyb prum"
$ws.Range("G8").Value = "var"

# New row 8: sample variance (VAR.S) of each question's answers, mirroring
# the AVERAGE row above it. H8 gets its own formula; I8:M8 share one.
$ws.Range("H8").Formula = "=VAR.S(H2:H6)"
$ws.Range("I8:M8").Formula = "=VAR.S(I2:I6)"

# Selection / active sheet now sit on the new variance row, and the
# "Vybrané otázky na parametry" sheet becomes the active tab.
[void]$ws.Range("H8:M8").Select()
$ws.Activate()

# Page setup gains an explicit paper size / orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Applied Albert's updates to 'Vybrané otázky na parametry'"
